$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: RowNum, A(Discount_ID), B(Shopkeeper_ID), C(Order_ID), D(Discount_Amount), E(Applied_By), F(Discount_Date), G(IsDeleted)
$rows = @(
    @(9,  8,  20, 8,  10, "System", "2025-03-28 18:18:25", 0),
    @(10, 9,  20, 9,  5,  "System", "2025-03-28 18:41:45", 0),
    @(11, 10, 16, 10, 0,  "System", "2025-03-28 19:07:07", 0),
    @(12, 11, 16, 11, 0,  "System", "2025-03-28 19:17:33", 0),
    @(13, 12, 16, 12, 0,  "System", "2025-03-28 19:18:43", 0),
    @(14, 13, 16, 13, 1,  "System", "2025-03-28 19:19:42", 0),
    @(15, 14, 16, 14, 1,  "System", "2025-03-28 19:21:47", 0),
    @(16, 15, 16, 15, 4,  "System", "2025-03-28 19:22:09", 0)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}
